# feat: add 2022-Q1 data
#
# - Renames the existing "总计" sheet to "2022-Q1" and replaces its content
#   with the new fund-holding breakdown for 2022-Q1.
# - Adds a brand-new "总计" sheet (placed after "2022-Q1", i.e. at the end)
#   containing the historical per-quarter summary, now including the new
#   2022-Q1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 0. Handy references / constants
# ---------------------------------------------------------------------
$xlPasteFormats = -4122

# A sheet that already carries the "header / index column" style (bold,
# centered, thin-box-bordered) we want to reuse so no redundant styles get
# minted.
$styleSourceSheet = $wb.Worksheets.Item("2021-Q4")
$headerStyleCell  = $styleSourceSheet.Cells.Item(1, 2)   # B1 -> s="2"

function Copy-HeaderStyle($cell) {
    $headerStyleCell.Copy() | Out-Null
    $cell.PasteSpecial($xlPasteFormats) | Out-Null
}

# Sets a cell to a literal TEXT value even when it looks like a number
# (mirrors typing `'12.3` into Excel -> stored as text, not coerced).
function Set-TextValue($cell, [string]$text) {
    $cell.Value = "'" + $text
}

# ---------------------------------------------------------------------
# 1. "总计" (sheetId 6) becomes "2022-Q1" with the new fund table
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$q1Header = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $cell = $q1.Cells.Item(1, $col)
    Copy-HeaderStyle $cell
    $cell.Value = $q1Header[$col - 2]
}

$q1Rows = @(
    @("004856", "广发中证全指建筑材料指数A", "13.72", "94.61", "2.44", "0.3348", 10),
    @("004857", "广发中证全指建筑材料指数C", "6.05",  "94.61", "2.44", "0.1476", 10),
    @("159610", "景顺长城中证500增强策略ETF", "8.45",  "98.35", "1.17", "0.0989", 8),
    @("159745", "国泰中证全指建筑材料交易型开放式指数证券投资基金", "3.76", "98.37", "2.58", "0.0970", 10),
    @("008672", "宝盈祥泽混合A", "4.71", "23.71", "1.18", "0.0556", 9),
    @("008673", "宝盈祥泽混合C", "3.59", "23.71", "1.18", "0.0424", 9),
    @("512590", "浦银安盛中证高股息精选ETF", "0.59", "96.43", "3.04", "0.0179", 2),
    @("516750", "富国中证全指建筑材料ETF", "0.47", "98.22", "2.54", "0.0119", 10),
    @("167702", "德邦量化优选股票(LOF)A", "0.54", "83.48", "2.05", "0.0111", 7),
    @("167703", "德邦量化优选股票(LOF)C", "0.41", "83.48", "2.05", "0.0084", 7),
    @("006143", "恒生前海中证质量成长低波动指数A", "0.06", "94.34", "3.02", "0.0018", 4),
    @("006144", "恒生前海中证质量成长低波动指数C", "0.01", "94.34", "3.02", "0.0003", 4),
    @("005770", "信达澳银中证沪港深高股息精选指数", "0.01", "92.47", "2.35", "0.0002", 8)
)

for ($i = 0; $i -lt $q1Rows.Count; $i++) {
    $r = $i + 2
    $row = $q1Rows[$i]

    $idxCell = $q1.Cells.Item($r, 1)
    Copy-HeaderStyle $idxCell
    $idxCell.Value = $i

    Set-TextValue $q1.Cells.Item($r, 2) $row[0]              # 基金代码
    $q1.Cells.Item($r, 3).Value = [string]$row[1]            # 基金名称
    Set-TextValue $q1.Cells.Item($r, 4) $row[2]             # 基金规模
    Set-TextValue $q1.Cells.Item($r, 5) $row[3]             # 股票总仓位
    Set-TextValue $q1.Cells.Item($r, 6) $row[4]             # 仓位占比
    Set-TextValue $q1.Cells.Item($r, 7) $row[5]             # 持有市值(亿元)
    $q1.Cells.Item($r, 8).Value = $row[6]                    # 仓位排名
}

# ---------------------------------------------------------------------
# 2. New "总计" sheet at the end, with the updated quarterly roll-up
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$totalHeader = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($col = 2; $col -le 4; $col++) {
    $cell = $total.Cells.Item(1, $col)
    Copy-HeaderStyle $cell
    $cell.Value = $totalHeader[$col - 2]
}

$totalRows = @(
    @("2022-Q1", 13, 0.83),
    @("2021-Q4", 4,  0.01),
    @("2021-Q3", 10, 0.28),
    @("2021-Q2", 4,  0.03),
    @("2021-Q1", 11, 0.18),
    @("2020-Q4", 9,  0.18)
)

for ($i = 0; $i -lt $totalRows.Count; $i++) {
    $r = $i + 2
    $row = $totalRows[$i]

    $idxCell = $total.Cells.Item($r, 1)
    Copy-HeaderStyle $idxCell
    $idxCell.Value = $i

    $total.Cells.Item($r, 2).Value = [string]$row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
}
